$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F binding header, mirrors the existing "Loan from $loan" style (D6)
# but as a plain (non-bordered) cell carrying the same pink row fill.
$ws.Range("F6").Value = 'Address from $address'
$ws.Range("F6").Interior.Color = 14145535

# The address binding used to be embedded in the condition text itself;
# now that it lives in F6, simplify the condition expression.
$ws.Range("F7").Value = 'state in ($param)'

$null = $ws.Range("F13").Select()
